$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the real ontology IRIs for the four rows that previously had the
#     placeholder "OBI or OBIB?" text (row 2 gets a brand-new D cell).
#     Written in reverse row order so the new shared-string entries land in
#     the same order as the target file (EUPATH_0000129 .. EUPATH_0000126).
$ws.Range("D6").Value = "http://purl.obolibrary.org/obo/EUPATH_0000129"
$ws.Range("D5").Value = "http://purl.obolibrary.org/obo/EUPATH_0000128"
$ws.Range("D3").Value = "http://purl.obolibrary.org/obo/EUPATH_0000127"
$ws.Range("D2").Value = "http://purl.obolibrary.org/obo/EUPATH_0000126"

# --- Update the view: scroll so column C is left-most and select D12.
$ws.Activate()
try { $excel.ActiveWindow.ScrollColumn = 3 } catch {}
$ws.Range("D12").Select()
